# Re-upload refresh of the "repayment_20250901_20250912" report:
#  - bump the report revision suffix in the sheet name: (2) -> (3)
#  - refresh Talk_time (H) figures for nearly every collector
#  - refresh Yandi Nugraha's (row 5) repayment count/amount and pending-amount-recovery rate
#  - refresh Fadilah Damayanti's (row 13) repayment count/amount and pending-amount-recovery rate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: Yandi Nugraha ---
$ws.Range("D5").Value = 28

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "16,534,223.00"
$ws.Range("E5").Style = "Normal"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5.97"
$ws.Range("G5").Style = "Normal"

$ws.Range("H5").Value = 16.088000000000001

# --- Row 13: Fadilah Damayanti ---
$ws.Range("D13").Value = 34

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "23,465,251.00"
$ws.Range("E13").Style = "Normal"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "7.69"
$ws.Range("G13").Style = "Normal"

$ws.Range("H13").Value = 11.930999999999999

# --- Talk_time (column H) refresh for the remaining collectors ---
$ws.Range("H2").Value = 12.566000000000001
$ws.Range("H3").Value = 17.14
$ws.Range("H4").Value = 9.89
$ws.Range("H6").Value = 12.077
$ws.Range("H7").Value = 8.3460000000000001
$ws.Range("H8").Value = 17.161999999999999
$ws.Range("H9").Value = 8.8659999999999997
$ws.Range("H10").Value = 9.1370000000000005
$ws.Range("H11").Value = 10.8
$ws.Range("H12").Value = 14.705
$ws.Range("H14").Value = 8.1579999999999995
$ws.Range("H15").Value = 6.984
$ws.Range("H16").Value = 7.2949999999999999
$ws.Range("H17").Value = 16.492999999999999
$ws.Range("H18").Value = 5.8010000000000002

# --- Sheet name bump (2) -> (3) ---
$ws.Name = "repayment_20250901_20250912 (3)"
